# Simulated Wild Card round and logged it
# Update the Home-row (row 2) target-depth stats on both the OFF and DEF
# sheets with the results of the simulated game.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 418
$wsOff.Range("C2").Value = 292
$wsOff.Range("D2").Value = 107
$wsOff.Range("E2").Value = 46
$wsOff.Range("G2").Value = 9

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 556
$wsDef.Range("C2").Value = 387
$wsDef.Range("D2").Value = 139
$wsDef.Range("E2").Value = 59
$wsDef.Range("G2").Value = 7
